$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
# "T" -> "Tobs", a new "variant" column is inserted after "d" (shifting the
# rest of the headers one column to the right), and a new trailing
# "avg_K_selective_cond" column is appended at T1.
$ws.Range("A1").Value = "N"
$ws.Range("B1").Value = "Tobs"
$ws.Range("C1").Value = "d"
$ws.Range("D1").Value = "variant"
$ws.Range("E1").Value = "overall"
$ws.Range("F1").Value = "overall_cond"
$ws.Range("G1").Value = "hom"
$ws.Range("H1").Value = "hom_cond"
$ws.Range("I1").Value = "rand_split"
$ws.Range("J1").Value = "rand_split_cond"
$ws.Range("K1").Value = "rand_selective"
$ws.Range("L1").Value = "rand_selective_cond"
$ws.Range("M1").Value = "recovery_split"
$ws.Range("N1").Value = "recovery_split_cond"
$ws.Range("O1").Value = "recovery_selective"
$ws.Range("P1").Value = "recovery_selective_cond"
$ws.Range("Q1").Value = "avg_K_split"
$ws.Range("R1").Value = "avg_K_split_cond"
$ws.Range("S1").Value = "avg_K_selective"
$ws.Range("T1").Value = "avg_K_selective_cond"

# --- Row 2 (existing data row, N=80, Tobs=20) ---------------------------
$ws.Range("A2").Value = 80
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = "overall_holds"
$ws.Range("E2").Value = 0.06
$ws.Range("F2").Value = 0.0625
$ws.Range("G2").Value = 0.0525
$ws.Range("H2").Value = 0.0225
$ws.Range("I2").Value = 0.0168718663913452
$ws.Range("J2").Value = 0.0196191134013549
$ws.Range("K2").Value = 0.0170813691750611
$ws.Range("L2").Value = 0.00920390130374723
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 3.485
$ws.Range("R2").Value = 3.9425
$ws.Range("S2").Value = 2.91
$ws.Range("T2").Value = 2.6575

# --- Row 3 (new data row, N=80, Tobs=50) --------------------------------
$ws.Range("A3").Value = 80
$ws.Range("B3").Value = 50
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = "overall_holds"
$ws.Range("E3").Value = 0.0525
$ws.Range("F3").Value = 0.045
$ws.Range("G3").Value = 0.0475
$ws.Range("H3").Value = 0.0575
$ws.Range("I3").Value = 0.0183582327084883
$ws.Range("J3").Value = 0.0132921262720714
$ws.Range("K3").Value = 0.0137318146155423
$ws.Range("L3").Value = 0.00711324528978235
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 3.0575
$ws.Range("R3").Value = 3.1675
$ws.Range("S3").Value = 2.6925
$ws.Range("T3").Value = 2.3
